# Adjusted MPCs (marginal propensity to consume) inputs and re-ran the FIM
# (Fiscal Impact Model), which refreshes both the "current" quarterly series
# and the "difference" (current - previous) series on Sheet 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Federal Contribution (current)
$ws.Range("N5").Value = -0.239

# Row 9: Federal Other Direct Aid Arp Contribution (current)
$ws.Range("K9").Value = -0.0422
$ws.Range("L9").Value = -0.0273
$ws.Range("M9").Value = -0.0277
$ws.Range("N9").Value = -0.0372
$ws.Range("O9").Value = -0.0357
$ws.Range("P9").Value = -0.0344
$ws.Range("Q9").Value = -0.0224
$ws.Range("R9").Value = -0.0152
$ws.Range("S9").Value = -0.0137
$ws.Range("T9").Value = -0.0138
$ws.Range("U9").Value = -0.0129
$ws.Range("V9").Value = -0.0143

# Row 14: Federal Subsidies Contribution (current)
$ws.Range("C14").Value = 0.3559
$ws.Range("D14").Value = 0.1922
$ws.Range("E14").Value = 0.0208
$ws.Range("F14").Value = -0.0658
$ws.Range("G14").Value = -0.1684
$ws.Range("H14").Value = -0.1661
$ws.Range("I14").Value = -0.0803
$ws.Range("J14").Value = -0.063
$ws.Range("K14").Value = -0.3297
$ws.Range("L14").Value = -0.3425
$ws.Range("M14").Value = -0.1399
$ws.Range("N14").Value = -0.1031
$ws.Range("O14").Value = -0.2902
$ws.Range("P14").Value = -0.2613
$ws.Range("Q14").Value = -0.1046
$ws.Range("R14").Value = -0.0568
$ws.Range("S14").Value = -0.0915
$ws.Range("T14").Value = -0.0676
$ws.Range("U14").Value = -0.0297
$ws.Range("V14").Value = -0.0777

# Row 16: Fiscal Impact (current)
$ws.Range("C16").Value = -2.452
$ws.Range("D16").Value = -2.9847
$ws.Range("E16").Value = -3.1702
$ws.Range("F16").Value = -3.7693
$ws.Range("G16").Value = -4.585
$ws.Range("H16").Value = -2.2495
$ws.Range("I16").Value = -0.5037
$ws.Range("J16").Value = 0.1538
$ws.Range("K16").Value = -0.1066
$ws.Range("L16").Value = 0.6392
$ws.Range("M16").Value = 0.1675
$ws.Range("N16").Value = -0.1075
$ws.Range("O16").Value = -0.6267
$ws.Range("P16").Value = -0.7305
$ws.Range("Q16").Value = -0.4342
$ws.Range("R16").Value = -0.698
$ws.Range("S16").Value = -0.6044
$ws.Range("T16").Value = -0.2413
$ws.Range("U16").Value = -0.0447
$ws.Range("V16").Value = -50.1495

# Row 18: Grants Contribution (current)
$ws.Range("N18").Value = -0.1184

# Row 20: Rebate Checks Arp Contribution (current)
$ws.Range("N20").Value = -0.1125
$ws.Range("O20").Value = -0.1052
$ws.Range("P20").Value = -0.102
$ws.Range("Q20").Value = -0.0999
$ws.Range("R20").Value = -0.0968
$ws.Range("S20").Value = -0.0938
$ws.Range("T20").Value = -0.001
$ws.Range("U20").Value = -0.001
$ws.Range("V20").Value = -0.0009

# Row 22: State Contribution (current)
$ws.Range("N22").Value = 0.1147

# Row 26: State Purchases Deflator Growth (current)
$ws.Range("N26").Value = 0.0074

# Row 28: State Subsidies Contribution (current)
$ws.Range("F28").Value = -0.0006
$ws.Range("G28").Value = -0.0007
$ws.Range("J28").Value = -0.0003
$ws.Range("K28").Value = -0.0009
$ws.Range("N28").Value = -0.0011
$ws.Range("O28").Value = -0.0025
$ws.Range("P28").Value = -0.0005
$ws.Range("R28").Value = -0.0006
$ws.Range("S28").Value = -0.0013
$ws.Range("T28").Value = -0.0003

# Row 33: Federal Contribution (difference)
$ws.Range("N33").Value = 0.1364

# Row 37: Federal Other Direct Aid Arp Contribution (difference)
$ws.Range("K37").Value = 0.06
$ws.Range("L37").Value = 0.009
$ws.Range("M37").Value = -0.0041
$ws.Range("N37").Value = -0.0112
$ws.Range("O37").Value = -0.0154
$ws.Range("P37").Value = -0.013
$ws.Range("Q37").Value = 0.0072
$ws.Range("R37").Value = 0.0055
$ws.Range("S37").Value = 0.0021
$ws.Range("T37").Value = -0.0031
$ws.Range("U37").Value = -0.011
$ws.Range("V37").Value = -0.0112

# Row 42: Federal Subsidies Contribution (difference)
$ws.Range("C42").Value = -0.0366
$ws.Range("D42").Value = -0.0427
$ws.Range("E42").Value = -0.0112
$ws.Range("F42").Value = -0.0099
$ws.Range("G42").Value = -0.0875
$ws.Range("H42").Value = -0.0831
$ws.Range("I42").Value = -0.0295
$ws.Range("J42").Value = -0.0177
$ws.Range("K42").Value = 0.1753
$ws.Range("L42").Value = 0.2099
$ws.Range("M42").Value = 0.0845
$ws.Range("N42").Value = 0.0728
$ws.Range("O42").Value = 0.0261
$ws.Range("P42").Value = -0.0243
$ws.Range("Q42").Value = 0.0001
$ws.Range("R42").Value = -0.0213
$ws.Range("S42").Value = -0.0626
$ws.Range("T42").Value = -0.0459
$ws.Range("U42").Value = -0.0167
$ws.Range("V42").Value = -0.0097

# Row 44: Fiscal Impact (difference)
$ws.Range("C44").Value = -0.0366
$ws.Range("D44").Value = -0.0427
$ws.Range("E44").Value = -0.0112
$ws.Range("F44").Value = -0.0099
$ws.Range("G44").Value = -0.0878
$ws.Range("H44").Value = -0.0831
$ws.Range("I44").Value = -0.0295
$ws.Range("J44").Value = -0.0178
$ws.Range("K44").Value = 0.2348
$ws.Range("L44").Value = 0.2189
$ws.Range("M44").Value = 0.4469
$ws.Range("N44").Value = 0.5838
$ws.Range("O44").Value = -0.0821
$ws.Range("P44").Value = -0.1187
$ws.Range("Q44").Value = -0.0787
$ws.Range("R44").Value = -0.1434
$ws.Range("S44").Value = -0.1737
$ws.Range("T44").Value = -0.0696
$ws.Range("U44").Value = -0.0423
$ws.Range("V44").Value = 0.0807

# Row 46: Grants Contribution (difference)
$ws.Range("N46").Value = 0.1371

# Row 48: Rebate Checks Arp Contribution (difference)
$ws.Range("N48").Value = 0.4836
$ws.Range("O48").Value = -0.105
$ws.Range("P48").Value = -0.1019
$ws.Range("Q48").Value = -0.0939
$ws.Range("R48").Value = -0.0968
$ws.Range("S48").Value = -0.0938
$ws.Range("T48").Value = -0.001
$ws.Range("U48").Value = -0.001
$ws.Range("V48").Value = -0.0009

# Row 50: State Contribution (difference)
$ws.Range("N50").Value = -0.1015

# Row 54: State Purchases Deflator Growth (difference)
$ws.Range("N54").Value = -0.0008

# Row 56: State Subsidies Contribution (difference)
$ws.Range("F56").Value = -0.0001
$ws.Range("G56").Value = -0.0003
$ws.Range("J56").Value = -0.0001
$ws.Range("K56").Value = -0.0005
$ws.Range("N56").Value = 0.0004
$ws.Range("O56").Value = 0.0018
$ws.Range("R56").Value = -0.0002
$ws.Range("S56").Value = -0.001
$ws.Range("V56").Value = 0.0001
